$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativação: 01/01/2018 -> 01/01/2021 -----------------------------------
# A plain ".Value = '01/01/2021'" assignment gets auto-detected by Excel as
# a date and turned into a serial number with a date number format. The
# source workbook stores this as a literal text string, so we round-trip it
# through a text formula and then convert that formula to a static value
# in place (PasteSpecial values-only). That keeps the cell a plain text
# (shared-string) cell with its original "General" number format/style.
$dateText = "01/01/2021"

$cell = $ws.Range("B8")
$cell.Formula = '="' + $dateText + '"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues

$cell = $ws.Range("C8")
$cell.Formula = '="' + $dateText + '"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues

$ws.Application.CutCopyMode = $false

# --- Docentes responsáveis --------------------------------------------------
$professor = "11079086 - Herlandí de Souza Andrade"
$ws.Range("B13").Value = $professor
$ws.Range("C13").Value = $professor

# --- Método: ----------------------------------------------------------------
$metodo = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Critério: ----------------------------------------------------------------
$criterio = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Norma de recuperação: ---------------------------------------------------
$norma = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Bibliografia: ------------------------------------------------------------
$bibliografia = "KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.KOTLER, P.; KARTAJAYA, H.; SETIAWAN, I. Marketing 4.0: do Tradicional ao Digital. São Paulo: Sextante, 2017.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L.  Marketing Essencial. 5 ed. São Paulo: Pearson, 2013.SANDHUSEN, R. L. Marketing Básico - Série Essencial. 3 ed. São Paulo: Saraiva, 2010.SAPIRO, Arão., CHIAVENATO, I. Planejamento Estratégico. Campus, 2ª. edição, 2010 KOTLER, P. Administração de Marketing, edição do milênio, revisão técnica de Prof. Arão Sapiro. Prentice-Hall, 2000. HOOLEY, Graham J.; PIERCY, Nigel F.; SAUNDERS, John A. Estratégia de Marketing e Posicionamento Competitivo tradução e revisão técnica: Prof. Arão Sapiro. Pearson Education do Brasil, 2001. SAPIRO, ARAO; GANGANA, MAURÍCIO; LIMA, MIGUEL; VILHENA, JOÃO BAPTISTA. Gestão de Marketing . FGV Editora, 2004. BOONE, L. e KURTZ, D.L. Marketing contemporâneo. 8ª ed. São Paulo, Livros Técnicos e Científicos, 1998. KOTLER, P; JATURISPITAK, S. e MAESINCIE, S. O marketing das nações. São Paulo, Futura, 1997. MARTINS, J.R. e BLECHER, N. O império das marcas. 2ª ed. São Paulo, Negócio Editora, 1997 THUROW, L.C. O futuro do capitalismo. 2ª ed. São Paulo, Rocco, 1997. VAZ, G. N. Marketing institucional. São Paulo, Pioneira, 1995. Bibliografia Complementar Artigos das Revistas: Marketing, Meio e Mensagem, Exame, Dinheiro, Revista da Escola de Administração da FEA-USP, Revista ESPM."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
